$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-01-20 Tuesday" "2026-01-21 Wednesday"

Replace-Text "99×61=6039" "93×70=6510"
Replace-Text "72×79=5688" "77×50=3850"
Replace-Text "63×13=819" "68×70=4760"
Replace-Text "65×76=4940" "62×93=5766"
Replace-Text "49×64=3136" "49×97=4753"
Replace-Text "58×64=3712" "71×35=2485"
Replace-Text "20×69=1380" "83×23=1909"
Replace-Text "75×79=5925" "51×87=4437"
Replace-Text "91×96=8736" "99×53=5247"
Replace-Text "99×73=7227" "17×44=748"
Replace-Text "24×60=1440" "33×77=2541"
Replace-Text "88×55=4840" "69×73=5037"
Replace-Text "60×46=2760" "37×73=2701"
Replace-Text "90×75=6750" "46×43=1978"
Replace-Text "54×33=1782" "50×76=3800"
Replace-Text "60×88=5280" "54×80=4320"
Replace-Text "12×83=996" "46×67=3082"
Replace-Text "48×56=2688" "15×76=1140"
Replace-Text "71×60=4260" "50×16=800"
Replace-Text "42×89=3738" "63×20=1260"
Replace-Text "62×35=2170" "31×20=620"
Replace-Text "99×51=5049" "53×71=3763"
Replace-Text "32×29=928" "58×82=4756"
Replace-Text "13×94=1222" "13×57=741"
Replace-Text "52×11=572" "86×56=4816"
